$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.048.68"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.77"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.75"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6213"
$ws.Range("E6").Value = "  -6.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.41"
$ws.Range("E8").Value = "  +5.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07370"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2923"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.71"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07673"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.828.66"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.963"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6629"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.08"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009068"
$ws.Range("E17").Value = "  +8.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.029"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.049.44"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.077.60"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "225.52"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.37"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.148"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.9999"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.61"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.422"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1357"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.80"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.493"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.038"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.054"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.199"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05243"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.845"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7321"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.289.82"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.749"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.303"
$ws.Range("E42").Value = "  +5.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9014"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.89"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.975.50"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5113"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.93"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.716"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3966"
$ws.Range("E51").Value = "  -1.78%  "
